# Scheduled runner update: refresh market-board derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# Leve profit tables on each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2431.585
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 2459.1155
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 7377.3465
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -7713.3465

$ws.Range("H100").Value = 1082.8422
$ws.Range("I100").Value = 864.2143
$ws.Range("J100").Value = 1695
$ws.Range("K100").Value = 864.2143
$ws.Range("L100").Value = 1695
$ws.Range("M100").Value = -323.2143
$ws.Range("N100").Value = -2777

$ws.Range("H106").Value = 1204.5
$ws.Range("I106").Value = 1018.125
$ws.Range("K106").Value = 1018.125
$ws.Range("M106").Value = -387.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 270
$ws.Range("I4").Value = 270
$ws.Range("K4").Value = 270
$ws.Range("M4").Value = -154

$ws.Range("H5").Value = 422.75
$ws.Range("I5").Value = 434
$ws.Range("J5").Value = 389
$ws.Range("K5").Value = 434
$ws.Range("L5").Value = 389
$ws.Range("M5").Value = -322
$ws.Range("N5").Value = -613

$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

$ws.Range("H45").Value = 1661.7428
$ws.Range("I45").Value = 1615.8462
$ws.Range("J45").Value = 1794.3334
$ws.Range("K45").Value = 1615.8462
$ws.Range("L45").Value = 1794.3334
$ws.Range("M45").Value = -1238.8462
$ws.Range("N45").Value = -2548.3334

$ws.Range("H63").Value = 3632.5
$ws.Range("I63").Value = 3361
$ws.Range("K63").Value = 3361
$ws.Range("M63").Value = -2675

$ws.Range("H66").Value = 3632.5
$ws.Range("I66").Value = 3361
$ws.Range("K66").Value = 16805
$ws.Range("M66").Value = -13373

$ws.Range("H97").Value = 1356.6666
$ws.Range("I97").Value = 1059.1666
$ws.Range("J97").Value = 1951.6666
$ws.Range("K97").Value = 1059.1666
$ws.Range("L97").Value = 1951.6666
$ws.Range("M97").Value = -563.1666
$ws.Range("N97").Value = -2943.6666

$ws.Range("H122").Value = 25002362
$ws.Range("I122").Value = 2937.3333
$ws.Range("K122").Value = 8811.999899999999
$ws.Range("M122").Value = -6361.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 422.75
$ws.Range("I4").Value = 434
$ws.Range("J4").Value = 389
$ws.Range("K4").Value = 434
$ws.Range("L4").Value = 389
$ws.Range("M4").Value = -319
$ws.Range("N4").Value = -619

$ws.Range("H94").Value = 1220.4
$ws.Range("I94").Value = 992.7143
$ws.Range("J94").Value = 1751.6666
$ws.Range("K94").Value = 992.7143
$ws.Range("L94").Value = 1751.6666
$ws.Range("M94").Value = -541.7143
$ws.Range("N94").Value = -2653.6666

$ws.Range("H99").Value = 1730.2941
$ws.Range("I99").Value = 1700.5555
$ws.Range("J99").Value = 1763.75
$ws.Range("K99").Value = 1700.5555
$ws.Range("L99").Value = 1763.75
$ws.Range("M99").Value = -202.5554999999999
$ws.Range("N99").Value = -4759.75

$ws.Range("H107").Value = 555.9474
$ws.Range("I107").Value = 462.14285
$ws.Range("J107").Value = 818.6
$ws.Range("K107").Value = 462.14285
$ws.Range("L107").Value = 818.6
$ws.Range("M107").Value = 1457.85715
$ws.Range("N107").Value = -4658.6

$ws.Range("H134").Value = 17219.293
$ws.Range("I134").Value = 1701.5106
$ws.Range("J134").Value = 57737.945
$ws.Range("K134").Value = 5104.531800000001
$ws.Range("L134").Value = 173213.835
$ws.Range("M134").Value = -2569.531800000001
$ws.Range("N134").Value = -178283.835

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 27504900
$ws.Range("I6").Value = 55000000
$ws.Range("J6").Value = 9800
$ws.Range("K6").Value = 55000000
$ws.Range("L6").Value = 9800
$ws.Range("M6").Value = -54999887
$ws.Range("N6").Value = -10026

$ws.Range("H22").Value = 180
$ws.Range("I22").Value = 172
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 172
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 178
$ws.Range("N22").Value = -900

$ws.Range("H58").Value = 1517107.4
$ws.Range("I58").Value = 2067407.9
$ws.Range("J58").Value = 3780.875
$ws.Range("K58").Value = 2067407.9
$ws.Range("L58").Value = 3780.875
$ws.Range("M58").Value = -2067204.9
$ws.Range("N58").Value = -4186.875

$ws.Range("H99").Value = 2579.6
$ws.Range("I99").Value = 2604
$ws.Range("J99").Value = 2569.1428
$ws.Range("K99").Value = 2604
$ws.Range("L99").Value = 2569.1428
$ws.Range("M99").Value = -1106
$ws.Range("N99").Value = -5565.1428

$ws.Range("H126").Value = 2579.6
$ws.Range("I126").Value = 2604
$ws.Range("J126").Value = 2569.1428
$ws.Range("K126").Value = 7812
$ws.Range("L126").Value = 7707.428400000001
$ws.Range("M126").Value = -5342
$ws.Range("N126").Value = -12647.4284

$ws.Range("H134").Value = 2524.516
$ws.Range("I134").Value = 1617.6666
$ws.Range("J134").Value = 3556.4482
$ws.Range("K134").Value = 4852.9998
$ws.Range("L134").Value = 10669.3446
$ws.Range("M134").Value = -2317.9998
$ws.Range("N134").Value = -15739.3446

$ws.Range("H136").Value = 1517107.4
$ws.Range("I136").Value = 2067407.9
$ws.Range("J136").Value = 3780.875
$ws.Range("K136").Value = 6202223.699999999
$ws.Range("L136").Value = 11342.625
$ws.Range("M136").Value = -6199673.699999999
$ws.Range("N136").Value = -16442.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.46875
$ws.Range("I2").Value = 8.5
$ws.Range("J2").Value = 41.791668
$ws.Range("K2").Value = 51
$ws.Range("L2").Value = 250.750008
$ws.Range("M2").Value = 62
$ws.Range("N2").Value = -476.750008

$ws.Range("H47").Value = 450
$ws.Range("I47").Value = 266.66666
$ws.Range("J47").Value = 1000
$ws.Range("K47").Value = 799.9999799999999
$ws.Range("L47").Value = 3000
$ws.Range("M47").Value = -368.9999799999999
$ws.Range("N47").Value = -3862

$ws.Range("H92").Value = 655.1539
$ws.Range("I92").Value = 608.4
$ws.Range("J92").Value = 684.375
$ws.Range("K92").Value = 1825.2
$ws.Range("L92").Value = 2053.125
$ws.Range("M92").Value = -577.1999999999998
$ws.Range("N92").Value = -4549.125

$ws.Range("H97").Value = 1973.4
$ws.Range("J97").Value = 2430.6667
$ws.Range("L97").Value = 7292.000100000001
$ws.Range("N97").Value = -8284.000100000001

$ws.Range("H132").Value = 1902.2084
$ws.Range("I132").Value = 2406.75
$ws.Range("J132").Value = 1649.9375
$ws.Range("K132").Value = 21660.75
$ws.Range("L132").Value = 14849.4375
$ws.Range("M132").Value = -19130.75
$ws.Range("N132").Value = -19909.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 65.666664
$ws.Range("I2").Value = 41.666668
$ws.Range("J2").Value = 89.666664
$ws.Range("K2").Value = 41.666668
$ws.Range("L2").Value = 89.666664
$ws.Range("M2").Value = 71.333332
$ws.Range("N2").Value = -315.666664

$ws.Range("H24").Value = 421052.6
$ws.Range("I24").Value = 5507500
$ws.Range("J24").Value = 14136.8
$ws.Range("K24").Value = 5507500
$ws.Range("L24").Value = 14136.8
$ws.Range("M24").Value = -5507327
$ws.Range("N24").Value = -14482.8

$ws.Range("H96").Value = 48950
$ws.Range("J96").Value = 48950
$ws.Range("L96").Value = 48950
$ws.Range("N96").Value = -54442

$ws.Range("H132").Value = 3849.8
$ws.Range("I132").Value = 1688.2
$ws.Range("J132").Value = 21142.6
$ws.Range("K132").Value = 5064.6
$ws.Range("L132").Value = 63427.8
$ws.Range("M132").Value = -2534.6
$ws.Range("N132").Value = -68487.79999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 994.75
$ws.Range("I22").Value = 1900
$ws.Range("J22").Value = 693
$ws.Range("K22").Value = 1900
$ws.Range("L22").Value = 693
$ws.Range("M22").Value = -1605
$ws.Range("N22").Value = -1283

$ws.Range("H27").Value = 994.75
$ws.Range("I27").Value = 1900
$ws.Range("J27").Value = 693
$ws.Range("K27").Value = 1900
$ws.Range("L27").Value = 693
$ws.Range("M27").Value = -1793
$ws.Range("N27").Value = -907

$ws.Range("H40").Value = 4811.875
$ws.Range("I40").Value = 4099
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 4099
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -3963
$ws.Range("N40").Value = -6272

$ws.Range("H100").Value = 4199.3887
$ws.Range("J100").Value = 6500.5
$ws.Range("L100").Value = 6500.5
$ws.Range("N100").Value = -7582.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 7526
$ws.Range("I3").Value = 5150
$ws.Range("J3").Value = 9902
$ws.Range("K3").Value = 5150
$ws.Range("L3").Value = 9902
$ws.Range("M3").Value = -5036
$ws.Range("N3").Value = -10130

$ws.Range("H14").Value = 25004166
$ws.Range("I14").Value = 50002500
$ws.Range("J14").Value = 5833
$ws.Range("K14").Value = 50002500
$ws.Range("L14").Value = 5833
$ws.Range("M14").Value = -50002332
$ws.Range("N14").Value = -6169

$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H122").Value = 4546.636
$ws.Range("I122").Value = 3312.875
$ws.Range("J122").Value = 7836.6665
$ws.Range("K122").Value = 9938.625
$ws.Range("L122").Value = 23509.9995
$ws.Range("M122").Value = -7488.625
$ws.Range("N122").Value = -28409.9995

$ws.Range("H132").Value = 1405.1224
$ws.Range("I132").Value = 711.7
$ws.Range("K132").Value = 2135.1
$ws.Range("M132").Value = 394.8999999999996
